$d = $word.ActiveDocument
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Comandos Fundamentais do GIt</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Versão do software:</w:t></w:r></w:p><w:p><w:r><w:t>Git –version</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Validade se possui algum repositório inicializado na pasta</w:t></w:r></w:p><w:p><w:r><w:t>Git status</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Criar um repositório:</w:t></w:r></w:p><w:p><w:r><w:t>Git init</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Adicionar </w:t></w:r><w:r><w:t xml:space="preserve">um </w:t></w:r><w:r><w:t>arquivo no repositório:</w:t></w:r></w:p><w:p><w:r><w:t>Git add “nome do arquivo”</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Adicionar vários arquivos no repositório:</w:t></w:r></w:p><w:p><w:r><w:t>Git add .</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Adicionar comentário </w:t></w:r><w:r><w:t xml:space="preserve">em um </w:t></w:r><w:r><w:t>arquivo adicionado ao repositório:</w:t></w:r></w:p><w:p><w:r><w:t>Git  commit</w:t></w:r><w:r><w:t xml:space="preserve"> “nome do arquivo”</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve"> -</w:t></w:r><w:r><w:t>m “</w:t></w:r><w:r><w:t>Estou enviando somente nome do arquivo</w:t></w:r><w:r><w:t>”</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Adicionar comentário único em todos arquivos que serão adicionados ao repositório:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Git commit </w:t></w:r><w:r><w:t xml:space="preserve">-a </w:t></w:r><w:r><w:t>-m “Estou enviando todos os arquivos”</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Enviar arquivos para servidor do Git:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Git Push </w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Receber mudanças que foram feitas no repositório:</w:t></w:r></w:p><w:p><w:r><w:t>Git pull</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Clonar o Repositório:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Git clone </w:t></w:r><w:r><w:t>https://github.com/lbsilva44/Curso_Github.git</w:t></w:r></w:p><w:p/><w:p/><w:p><w:r><w:t>Git branch -m main</w:t></w:r></w:p><w:p><w:r><w:t>Git branch -m master</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Git Push -u origin master  </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Content.InsertXML($xml)
Write-Host "Paragraphs after edit: " $d.Paragraphs.Count
